$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 14 de Octubre de 2020 a las 16:20"

# Data updates per row: B=Casos totales, C=Nuevos casos, D=Casos activos,
# E=Recuperados, G=Muertes hoy, H=Muertes
$rows = @(
    @{ R=4;   B=8098612; C=8359; D=5229233; E=2648379; G=127; H=221000 },
    @{ R=18;  B=413215;  C=3857; D=347396;  E=55798;   G=51;  H=10021 },
    @{ R=25;  B=337472;  C=1793; D=281900;  E=45826;   G=6;   H=9746 },
    @{ R=27;  B=297501;  C=849;  D=248940;  E=46494;   G=12;  H=2067 },
    @{ R=42;  B=113269;  C=532;  D=105236;  E=7357;    G=4;   H=676 },
    @{ R=57;  B=76272;   C=0;    D=72164;   E=3821;    G=2;   H=287 },
    @{ R=58;  B=68704;   C=2823; D=49800;   E=16796;   G=3;   H=2108 },
    @{ R=71;  B=45658;   C=458;  D=39304;   E=5961;    G=2;   H=393 },
    @{ R=77;  B=35251;   C=245;  D=31536;   E=2947;    G=1;   H=768 },
    @{ R=80;  B=31655;   C=482;  D=24074;   E=6623;    G=11;  H=958 },
    @{ R=95;  B=15888;   C=97;   D=11863;   E=3748;    G=0;   H=277 },
    @{ R=97;  B=15616;   C=29;   D=14850;   E=421;     G=0;   H=345 },
    @{ R=98;  B=15331;   C=24;   D=13571;   E=1444;    G=1;   H=316 },
    @{ R=107; B=10336;   C=39;   D=9245;    E=1012;    G=0;   H=79 },
    @{ R=110; B=10069;   C=124;  D=6531;    E=3443;    G=0;   H=95 },
    @{ R=145; B=3642;    C=6;    D=2595;    E=929;     G=0;   H=118 },
    @{ R=164; B=1833;    C=11;   D=961;     E=830;     G=0;   H=42 }
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
}
